$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "n_analysis(_string_)" -> "n_analysis(string)"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("n_analysis(_string_)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "n_analysis(string)", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 2: ", or `brackets'." -> ", or" + " " + "brackets"(styled) + "."
# ---------------------------------------------------------------------
$rightQuote = [char]0x2019
$oldTail = ", or ``brackets" + $rightQuote + "."
$fullText = $d.Content.Text
$tailIdx = $fullText.IndexOf($oldTail)

$tailRange = $d.Range($tailIdx, $tailIdx + $oldTail.Length)
$tailRange.Text = ", or "

$bracketsStart = $tailIdx + 5
$bracketsRange = $d.Range($bracketsStart, $bracketsStart)
$bracketsRange.InsertAfter("brackets")

$bracketsStyleRange = $d.Range($bracketsStart, $bracketsStart + 8)
$bracketsStyleRange.Style = "Verbatim Char"

$periodRange = $d.Range($bracketsStart + 8, $bracketsStart + 8)
$periodRange.InsertAfter(".")

# ---------------------------------------------------------------------
# Change 3a: remove the "When denominators or missing data summaries..."
# sentence (including the su_decimal(#)/miss_decimal(#) verbatim runs)
# from the paragraph, leaving it ending at "...the summaries."
# ---------------------------------------------------------------------
$fullText = $d.Content.Text
$keepMarker = "is specified the default is to place columns containing counts of nonmissing observations in each group before the columns containing the summaries."
$dropMarker = "can be used to independently control the number of decimal places reported for summary statistics and the percent of missing/nonmissing observations."

$keepIdx = $fullText.IndexOf($keepMarker)
$keepEnd = $keepIdx + $keepMarker.Length
$dropIdx = $fullText.IndexOf($dropMarker)
$dropEnd = $dropIdx + $dropMarker.Length

$removeRange = $d.Range($keepEnd, $dropEnd)
$removeRange.Text = ""

# ---------------------------------------------------------------------
# Change 3b: "N 2" -> "N 0" and "Summary 2" -> "Summary 0" in the
# "post `postname'" source-code line
# ---------------------------------------------------------------------
$d.Content.Find.Execute("(""N 2"") (""N Overall"") (""Summary 1"") (""Summary 2"")", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "(""N 0"") (""N Overall"") (""Summary 1"") (""Summary 0"")", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3c: drop "  miss_decimal(2) su_decimal(0)" from the
# "pt_base age" source-code line
# ---------------------------------------------------------------------
$d.Content.Find.Execute("cat_col  n_analysis(cols)  miss_decimal(2) su_decimal(0)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "cat_col  n_analysis(cols)", 2) | Out-Null

# ---------------------------------------------------------------------
# Change 3d: drop "  miss_decimal(2) decimal(1)" from the
# "pt_base qol" source-code line
# ---------------------------------------------------------------------
$d.Content.Find.Execute("cat_col  n_analysis(cols)  miss_decimal(2) decimal(1)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "cat_col  n_analysis(cols)", 2) | Out-Null
